$d = $word.ActiveDocument

# Locate the paragraph that ends the "Categories va About us" bullet -
# the new bullet about finishing the admin area needs to land right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Categories*About us*") {
        $target = $p
    }
}

# Grab the exact original OOXML for that paragraph (preserves its w14:paraId,
# rsid attributes, run formatting, etc. byte-for-byte) so re-inserting it
# alongside the new paragraph doesn't disturb it.
$full = $target.Range.WordOpenXML
$bodyIdx = $full.IndexOf("<w:body>")
$afterBody = $full.Substring($bodyIdx + 8)
$endIdx = $afterBody.IndexOf("</w:p>")
$origParXml = $afterBody.Substring(0, $endIdx + 6)

# The new paragraph to insert right after it.
$newParXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">-Hoàn thành cơ bản admin </w:t></w:r></w:p>'

$combined = $origParXml + $newParXml

# Collapsing the target paragraph's range to its end and inserting XML there
# replaces that paragraph's content - so we feed back its own (unchanged)
# markup followed by the new paragraph, which nets out to "insert after".
$r = $target.Range.Duplicate()
$r.Collapse(0)
$null = $r.InsertXML($combined)
